# "Remove deadzone in old manual"
#
# The last slide of the deck (slide 5) shows the analog-stick diagram from
# the old firmware manual. Over it sits a small annotation callout that is
# no longer accurate and is being removed:
#   - "左大括号 1"      (id 2)  - the left-brace connector bracketing the dial
#   - "文本框 8"        (id 9)  - the "Analog" / "Deadzone" text box
#   - "矩形: 圆角 9"    (id 10) - the "1.4°" degree pill
#   - "矩形: 圆角 17"   (id 18) - the "2.8°" degree pill
#   - "矩形: 圆角 20"   (id 21) - the "4.2°" degree pill
#
# The rest of the slide (including the other brace/label group, ids 12-17)
# is left untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

$namesToRemove = @("左大括号 1", "文本框 8", "矩形: 圆角 9", "矩形: 圆角 17", "矩形: 圆角 20")

foreach ($targetName in $namesToRemove) {
    for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
        $shp = $s.Shapes.Item($i)
        if ($shp.Name -eq $targetName) {
            $shp.Delete()
            break
        }
    }
}
